$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Grow the table by two rows (Table UR, Support UR) before touching any
# other cell so the table range / autofilter / dimension follow along.
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null

# Existing rows: values re-expressed in millimeters (x10 of previous cm values)
$ws.Range("B2").Value = 2000
$ws.Range("C2").Value = 1200
$ws.Range("D2").Value = 30

$ws.Range("B3").Value = 150
$ws.Range("C3").Value = 25
$ws.Range("D3").Value = 15

$ws.Range("B4").Value = 260
$ws.Range("C4").Value = 225
$ws.Range("D4").Value = 105

# New rows' labels (set before header renames below so shared strings land
# in the same order as the target workbook)
$ws.Range("A5").Value = "Table UR"
$ws.Range("A6").Value = "Support UR"

# Header renames: cm -> mm (updates both sheet cells and the table columns)
$ws.Range("D1").Value = "Z (mm)"
$ws.Range("C1").Value = "Y (mm)"
$ws.Range("B1").Value = "X (mm)"

# New rows' dimension values
$ws.Range("B5").Value = 870
$ws.Range("C5").Value = 580
$ws.Range("D5").Value = 28

$ws.Range("B6").Value = 500
$ws.Range("C6").Value = 500
$ws.Range("D6").Value = 40

# Match the final selection left by the author
$ws.Range("B6").Select()
